$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row: rename the Start/End Date headers to include the date format
# hint, matching the new explicit dd/mm/yyyy display format applied below.
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "Start Date (DD/MM/YYYY)"
$ws.Range("I1").Value = "End Date (DD/MM/YYYY)"

# ---------------------------------------------------------------------------
# Chapter 5 (row 6) is now finished: all 25 exercises done.
# ---------------------------------------------------------------------------
$ws.Range("D6").Value = 25

# ---------------------------------------------------------------------------
# Dates: replace the old TODAY()-based volatile formulas with a fixed start
# date plus a consistent "end = start + days - 1" chain, re-based on row 2's
# literal start date of 44413 (05 Aug 2021).
# ---------------------------------------------------------------------------
$ws.Range("H2").Value = 44413
$ws.Range("I2").Formula = "=H2+G2-1"

$ws.Range("H3").Formula = "=I2+1"
$ws.Range("I3").Formula = "=H3+G3-1"

$ws.Range("H4").Formula = "=I3+1"
$ws.Range("I4").Formula = "=H4+G4-1"

$ws.Range("H5").Formula = "=I4+1"
$ws.Range("I5").Formula = "=H5+G5-1"

$ws.Range("H6").Formula = "=I5+1"
$ws.Range("I6").Formula = "=H6+G6-1"

$ws.Range("H7").Formula = "=I6+1"
$ws.Range("I7").Formula = "=H7+G7-1"

$ws.Range("H8").Formula = "=I7+1"
$ws.Range("I8").Formula = "=H8+G8-1"

$ws.Range("H9").Formula = "=I8+1"
$ws.Range("I9").Formula = "=H9+G9-1"

$ws.Range("H10").Formula = "=I9+1"
$ws.Range("I10").Formula = "=H10+G10-1"

$ws.Range("H11").Formula = "=I10+1"
$ws.Range("I11").Formula = "=H11+G11-1"

$ws.Range("H12").Formula = "=I11+1"
$ws.Range("I12").Formula = "=H12+G12-1"

$ws.Range("H13").Formula = "=I12+1"
$ws.Range("I13").Formula = "=H13+G13-1"

$ws.Range("H14").Formula = "=I13+1"
$ws.Range("I14").Formula = "=H14+G14-1"

$ws.Range("H15").Formula = "=I14+1"
$ws.Range("I15").Formula = "=H15+G15-1"

$ws.Range("H16").Formula = "=I15+1"
$ws.Range("I16").Formula = "=H16+G16-1"

$ws.Range("H17").Formula = "=I16+1"
$ws.Range("I17").Formula = "=H17+G17-1"

$ws.Range("H18").Formula = "=I17+1"
$ws.Range("I18").Formula = "=H18+G18-1"

$ws.Range("H19").Formula = "=I18+1"
$ws.Range("I19").Formula = "=H19+G19-1"

$ws.Range("H20").Formula = "=I19+1"
$ws.Range("I20").Formula = "=H20+G20-1"

# Summary row: overall start/end date span.
$ws.Range("H21").Formula = "=H2"
$ws.Range("I21").Formula = "=I20"

# ---------------------------------------------------------------------------
# Apply the explicit dd/mm/yyyy custom date format to the Start/End Date
# data + summary rows (the header row stays plain text, no date format).
# ---------------------------------------------------------------------------
$ws.Range("H2:I21").NumberFormat = "dd/mm/yyyy;@"

# ---------------------------------------------------------------------------
# Selection, as last left by the author after editing row 6 (Chapter 5).
# ---------------------------------------------------------------------------
$ws.Range("G6:I6").Select()
